{"js": "// The template document has two near-identical blocks (one for\n// \"${consumables_block}\", one for \"${services_block}\"), each with a\n// label paragraph that (before this edit) reads\n// \"\u0417\u0430\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0440\u0430\u0441\u0445\u043e\u0434\u043d\u044b\u0435 \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u044b:\" right before its table. The diff\n// only touches the one that belongs to the services block, changing its\n// text to \"\u041e\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0443\u0441\u043b\u0443\u0433\u0438:\" (while keeping the run formatting\n// untouched). Disambiguate by walking the body paragraphs in order and\n// picking the label paragraph that follows the \"${services_block}\"\n// marker (i.e. the *last* match, since the services block comes second).\n\nconst OLD_TEXT = \"\u0417\u0430\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0440\u0430\u0441\u0445\u043e\u0434\u043d\u044b\u0435 \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u044b:\";\nconst NEW_TEXT_PART1 = \"\u041e\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0443\u0441\u043b\u0443\u0433\u0438\";\nconst NEW_TEXT_PART2 = \":\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === OLD_TEXT) {\n    // Keep walking \u2014 the services-block occurrence is the later one.\n    target = p;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the '\" + OLD_TEXT + \"' paragraph to update.\");\n}\n\n// Replace the run's text in place (preserves its run-level formatting:\n// Droid Serif, sz 21, b/i/caps/smallCaps=false, spacing=0) and then\n// append the trailing \":\" right after it, mirroring the diff's split\n// into two runs with identical rPr.\nconst range = target.getRange();\nconst firstRange = range.insertText(NEW_TEXT_PART1, \"Replace\");\nfirstRange.insertText(NEW_TEXT_PART2, \"After\");\n\nawait context.sync();\n", "ps1": "# The template document has two near-identical blocks (one for\n# \"${consumables_block}\", one for \"${services_block}\"), each with a label\n# paragraph that (before this edit) reads \"\u0417\u0430\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0440\u0430\u0441\u0445\u043e\u0434\u043d\u044b\u0435 \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u044b:\"\n# right before its table. The diff only touches the one belonging to the\n# services block, changing its text to \"\u041e\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0443\u0441\u043b\u0443\u0433\u0438:\" while leaving\n# the run formatting untouched. Disambiguate by scanning all paragraphs in\n# document order and keeping the LAST paragraph whose text equals the old\n# label (the services block comes after the consumables block), mirroring\n# the Office.js version's logic.\n\n$d = $word.ActiveDocument\n\n$oldText = \"\u0417\u0430\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0440\u0430\u0441\u0445\u043e\u0434\u043d\u044b\u0435 \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u044b:\"\n$newPart1 = \"\u041e\u043a\u0430\u0437\u0430\u043d\u043d\u044b\u0435 \u0443\u0441\u043b\u0443\u0433\u0438\"\n$newPart2 = \":\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq ($oldText + \"`r\")) {\n        $target = $p\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate the '$oldText' paragraph to update.\"\n}\n\n# Replace the run's text in place (preserves its run-level formatting:\n# Droid Serif, sz 21, b/i/caps/smallCaps=false, spacing=0), then append the\n# trailing \":\" right after it - mirrors the diff's split into two runs\n# sharing the same rPr.\n$rng = $target.Range\n$rng.MoveEnd(1, -1) | Out-Null\n$rng.Text = $newPart1\n\n$after = $rng.Duplicate\n$after.Collapse(0)\n$after.InsertAfter($newPart2)\n"}
